# Fruta / hortaliza, semanal
# Insert a new weekly record at row 195 (shifting the existing rows 195-205 down
# to 196-206), matching the weekly refresh of the "Vega Modelo de Temuco -
# Berenjena" consolidated dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row before the current row 195; this pushes the existing
# rows 195-205 down to 196-206 and keeps column D's date style (s="2").
$ws.Rows.Item(195).Insert()

# Populate the newly inserted row 195 with this week's record.
$ws.Cells.Item(195, 1).Value = 10
$ws.Cells.Item(195, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(195, 3).Value = "La Araucanía"
$ws.Cells.Item(195, 4).Value = 44516
$ws.Cells.Item(195, 5).Value = 9
$ws.Cells.Item(195, 6).Value = 100112001
$ws.Cells.Item(195, 7).Value = "Berenjena"
$ws.Cells.Item(195, 8).Value = "Sin especificar"
$ws.Cells.Item(195, 9).Value = "Primera"
$ws.Cells.Item(195, 10).Value = 15
$ws.Cells.Item(195, 11).Value = 10000
$ws.Cells.Item(195, 12).Value = 10000
$ws.Cells.Item(195, 13).Value = 10000
$ws.Cells.Item(195, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(195, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(195, 16).Value = 167
$ws.Cells.Item(195, 17).Value = 60
$ws.Cells.Item(195, 18).Value = "Hortaliza"
